$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new data row at row 552 (this pushes the existing row 552
# and everything below it down by one row, turning A1:R659 into A1:R660).
$ws.Rows.Item(552).Insert()

# Populate the newly inserted row 552 with the new record's data.
$ws.Cells.Item(552, 1).Value = 9
$ws.Cells.Item(552, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(552, 3).Value = "Metropolitana"
$ws.Cells.Item(552, 4).Value = 45258
$ws.Cells.Item(552, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(552, 5).Value = 13
$ws.Cells.Item(552, 6).Value = 100112039
$ws.Cells.Item(552, 7).Value = "Ciboulette"
$ws.Cells.Item(552, 8).Value = "Sin especificar"
$ws.Cells.Item(552, 9).Value = "Primera"
$ws.Cells.Item(552, 10).Value = 440
$ws.Cells.Item(552, 11).Value = 1000
$ws.Cells.Item(552, 12).Value = 1200
$ws.Cells.Item(552, 13).Value = 1064
$ws.Cells.Item(552, 14).Value = "`$/docena de atados"
$ws.Cells.Item(552, 15).Value = "Región Metropolitana"
$ws.Cells.Item(552, 16).Value = 355
$ws.Cells.Item(552, 17).Value = 3
$ws.Cells.Item(552, 18).Value = "Hortaliza"
